$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7135302424430847
$ws.Range("B1").Value = 0.58637934923172
$ws.Range("C1").Value = 0.5236685276031494
$ws.Range("D1").Value = 0.5788750648498535
$ws.Range("E1").Value = 0.7239269018173218
